$wb = $excel.ActiveWorkbook

# --- invalid_email sheet: unify the two validation-error messages into a single
#     shorter message, and move the active selection.
$ws = $wb.Worksheets.Item("invalid_email")
$ws.Range("C2").Value = "Please check username"
$ws.Range("C3").Value = "Please check username"
[void]$ws.Range("F13").Select()
